$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pollutant")
$tbl = $ws.ListObjects.Item("tbl_pollutant5")

# Insert a new row at worksheet row 37, shifting existing rows (37-45) down to (38-46)
$ws.Rows.Item(37).Insert()

$ws.Range("A37").Value = "NOx"
$ws.Range("B37").Value = "NO2"

# Grow the table range to include the newly inserted row
$tbl.Resize($ws.Range("A1:B46"))

$ws.Range("B37").Select()
